{"js": "// Office.js (Word JavaScript API) script implementing the\n// \"Code Inspection files added\" edit:\n//\n//   1. After the \"Code Inspection Meeting\" heading (and the blank paragraph\n//      that follows it) insert five narrative paragraphs describing the\n//      code inspection meeting, each with a 720-twip first-line indent.\n//   2. Tag the run of text \"Code\" (the later \"## Code\" heading) with a\n//      <w:lastRenderedPageBreak/> marker.\n//\n// Paragraphs are located by their trimmed text rather than fixed indices.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Part 1: insert the five \"Code Inspection Meeting\" narrative paragraphs\n// ---------------------------------------------------------------------\nlet headingIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Code Inspection Meeting\") {\n    headingIndex = i;\n    break;\n  }\n}\nif (headingIndex === -1) {\n  throw new Error(\"Could not find the 'Code Inspection Meeting' heading paragraph.\");\n}\n\nconst blankParagraph = paragraphs.items[headingIndex + 1];\n\nconst newParagraphsOoxml = \"<?xml version=\\\"1.0\\\" encoding=\\\"UTF-8\\\" standalone=\\\"yes\\\"?><pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\"><pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\"><pkg:xmlData><w:document xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\"><w:body><w:p><w:pPr><w:ind w:firstLine=\\\"720\\\"/></w:pPr><w:r><w:t xml:space=\\\"preserve\\\">We chose the class that implements most of the business logic for our system, to do the code inspection on it. The name of the class is Controller.java. The inspection took place in the Room 633 MUDD. The meeting started at around 11:45 am and finished around 12:20 pm. Since there was a lot of code in the \u201cunit\u201d, we could not go through all of the code in time; however, we got through the important functionality in time. </w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine=\\\"720\\\"/></w:pPr><w:r><w:t>Michael Glass was the appointed reader and Waseem Ilahi was the recorder. And obviously the TA (Jonathan Bell) was the moderator; however professor Kaiser was also present at the meeting.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine=\\\"720\\\"/></w:pPr><w:r><w:t>The logic of the code was all satisfactory and there was no problem with the code in that context. The TA went through the checklist to make sure each point was covered in the inspection. All the checkpoints were satisfied, except there was a little problem with the last two.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine=\\\"720\\\"/></w:pPr><w:r><w:t>We found a little inconsistency in the code, in terms of the usage of \u201ctrue\u201d and false. Also there was one other rather large Boolean expression, that might confuse the reader (suggested the TA). There was also some concern with the comments and there placement. The last thing that the TA suggested wa</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> to divide the main \u201cpart/unit\u201d of the inspected unit</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">, into logical components. </w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine=\\\"720\\\"/></w:pPr><w:r><w:t xml:space=\\\"preserve\\\">From the code files attached with folder containing this document, we can see that the modified \u201cunit\u201d implements the changes suggested. A \u201cdiff\u201d of the two sources will show the exact changes. </w:t></w:r><w:r><w:t>The main method (handle()) is left as a few calls to other methods, that act as a fork to take the control over to the appropriate method to handle the feature the user wants to use.</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> This process simplifies the main method and divides the entire logic into smaller pieces. </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> </w:t></w:r><w:r><w:t>As mentioned earlier, the changes can be seen in the</w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\"> file named \u201ccode_inspection_unit_modified.java\u201d.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n// Inserting at the \"End\" of the (empty) paragraph that follows the heading\n// appends the new paragraphs as siblings right after it, leaving the\n// paragraph that used to follow (the \"Defect Log\" heading) untouched.\nblankParagraph.insertOoxml(newParagraphsOoxml, Word.InsertLocation.end);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Part 2: tag the \"Code\" heading run with <w:lastRenderedPageBreak/>\n// ---------------------------------------------------------------------\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"text\");\nawait context.sync();\n\nlet codeIndex = -1;\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  if (paragraphs2.items[i].text.trim() === \"Code\") {\n    codeIndex = i;\n    break;\n  }\n}\nif (codeIndex === -1) {\n  throw new Error(\"Could not find the 'Code' heading paragraph.\");\n}\n\nconst codeParagraph = paragraphs2.items[codeIndex];\nconst pageBreakOoxml = \"<?xml version=\\\"1.0\\\" encoding=\\\"UTF-8\\\" standalone=\\\"yes\\\"?><pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\"><pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\"><pkg:xmlData><w:document xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>Code</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\n// \"Replace\" the whole paragraph's content with the same text plus the\n// <w:lastRenderedPageBreak/> marker prepended to the run.\ncodeParagraph.insertOoxml(pageBreakOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Word COM interop script: apply the \"Code Inspection files added\" edit.\n#\n# 1) After the \"Code Inspection Meeting\" heading (and its following blank\n#    paragraph) insert five narrative paragraphs, each indented with a\n#    first-line indent of 720 twips (0.5in).\n# 2) Mark the run of text \"Code\" (the \"## Code\" heading further down) with\n#    a <w:lastRenderedPageBreak/> element.\n#\n# Paragraphs are located by their exact trimmed text instead of hard-coded\n# indices so the script is resilient to minor structural differences.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Part 1: insert the five \"Code Inspection Meeting\" narrative paragraphs\n# ---------------------------------------------------------------------\n$headingIdx = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $para = $d.Paragraphs($i)\n    if ($para.Range.Text.Trim() -eq \"Code Inspection Meeting\") {\n        $headingIdx = $i\n        break\n    }\n}\n\nif ($headingIdx -eq -1) {\n    throw \"Could not find the 'Code Inspection Meeting' heading paragraph.\"\n}\n\n$blankPara = $d.Paragraphs($headingIdx + 1)\n$insertAt = $blankPara.Range.End\n$insertRange = $d.Range($insertAt, $insertAt)\n\n$newParasXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:ind w:firstLine=\"720\"/></w:pPr><w:r><w:t xml:space=\"preserve\">We chose the class that implements most of the business logic for our system, to do the code inspection on it. The name of the class is Controller.java. The inspection took place in the Room 633 MUDD. The meeting started at around 11:45 am and finished around 12:20 pm. Since there was a lot of code in the \u201cunit\u201d, we could not go through all of the code in time; however, we got through the important functionality in time. </w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine=\"720\"/></w:pPr><w:r><w:t>Michael Glass was the appointed reader and Waseem Ilahi was the recorder. And obviously the TA (Jonathan Bell) was the moderator; however professor Kaiser was also present at the meeting.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine=\"720\"/></w:pPr><w:r><w:t>The logic of the code was all satisfactory and there was no problem with the code in that context. The TA went through the checklist to make sure each point was covered in the inspection. All the checkpoints were satisfied, except there was a little problem with the last two.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine=\"720\"/></w:pPr><w:r><w:t>We found a little inconsistency in the code, in terms of the usage of \u201ctrue\u201d and false. Also there was one other rather large Boolean expression, that might confuse the reader (suggested the TA). There was also some concern with the comments and there placement. The last thing that the TA suggested wa</w:t></w:r><w:r><w:t>s</w:t></w:r><w:r><w:t xml:space=\"preserve\"> to divide the main \u201cpart/unit\u201d of the inspected unit</w:t></w:r><w:r><w:t xml:space=\"preserve\">, into logical components. </w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine=\"720\"/></w:pPr><w:r><w:t xml:space=\"preserve\">From the code files attached with folder containing this document, we can see that the modified \u201cunit\u201d implements the changes suggested. A \u201cdiff\u201d of the two sources will show the exact changes. </w:t></w:r><w:r><w:t>The main method (handle()) is left as a few calls to other methods, that act as a fork to take the control over to the appropriate method to handle the feature the user wants to use.</w:t></w:r><w:r><w:t xml:space=\"preserve\"> This process simplifies the main method and divides the entire logic into smaller pieces. </w:t></w:r><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:r><w:t>As mentioned earlier, the changes can be seen in the</w:t></w:r><w:r><w:t xml:space=\"preserve\"> file named \u201ccode_inspection_unit_modified.java\u201d.</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$insertRange.InsertXML($newParasXml)\n\n# InsertXML folds its LAST paragraph into whatever paragraph originally sat\n# at the insertion point (here: \"Defect Log\"). We padded the injected XML\n# with a throw-away trailing <w:p/> to keep \"Defect Log\" intact; now remove\n# that left-over empty paragraph that sits directly before \"Defect Log\".\n$count2 = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count2; $i++) {\n    $para = $d.Paragraphs($i)\n    if ($para.Range.Text.Trim() -eq \"Defect Log\") {\n        $prev = $d.Paragraphs($i - 1)\n        if ($prev.Range.Text.Trim() -eq \"\") {\n            $prev.Range.Delete()\n        }\n        break\n    }\n}\n\n# ---------------------------------------------------------------------\n# Part 2: tag the \"Code\" heading run with <w:lastRenderedPageBreak/>\n# ---------------------------------------------------------------------\n$count3 = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count3; $i++) {\n    $para = $d.Paragraphs($i)\n    if ($para.Range.Text.Trim() -eq \"Code\") {\n        $start = $para.Range.Start\n        $end = $para.Range.End - 1\n        $runRange = $d.Range($start, $end)\n        $pageBreakXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>Code</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n        $runRange.InsertXML($pageBreakXml)\n        break\n    }\n}\n"}
